$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the has_solution boolean flags (F column)
$ws.Range("F8").Value = $false
$ws.Range("F10").Value = $true
$ws.Range("F19").Value = $true
$ws.Range("F20").Value = $true
$ws.Range("F22").Value = $false
$ws.Range("F32").Value = $true
$ws.Range("F37").Value = $true
$ws.Range("F44").Value = $true

# Update the top-10 similar word lists (G column)
$ws.Range("G2").Value = "answer, response, liquid_bleach, bleach_liquor, gram's_solution, leading_question, question_of_law, rejoinder, spirits_of_ammonia, evasive_answer"
$ws.Range("G3").Value = "cry, crying, snivel, blue_murder, whimper, wailing, blubberer, lament, weeper, bawler"
$ws.Range("G4").Value = "deer, elk, mule_deer, antler, pere_david's_deer, fallow_deer, red_deer, cervus, roe_deer, wapiti"
$ws.Range("G5").Value = "floral_leaf, flower, flower_bud, umbrellawort, common_dandelion, easter_daisy, ray_flower, dandelion_green, tidytips, petal"
$ws.Range("G6").Value = "male_horse, filly, foal, stallion, ridgeling, horse, polo_pony, gelding, broodmare, pony"
$ws.Range("G8").Value = "pure_mathematics, vector_algebra, integral_calculus, quadratics, spherical_trigonometry, differential_calculus, matrix_algebra, ptyalith, calculus_of_variations, affine_geometry"
$ws.Range("G9").Value = "push, bellpull, thrust, nudge, tug, pedal_point, lever, squeeze, sustaining_pedal, brake_pedal"
$ws.Range("G10").Value = "right, dextrorotation, right_of_action, right_of_search, rightism"
$ws.Range("G11").Value = "influenza, contagious_disease, tumor_virus, upper_respiratory_infection, viral_infection, slow_virus, swine_influenza, respiratory_syncytial_virus, communicable_disease, asian_influenza"
$ws.Range("G12").Value = "celestial_body, outer_planet, superior_planet"
$ws.Range("G13").Value = "rough_fish, fish, food_fish, panfish, freshwater_fish, fishing, game_fish, scombroid, striped_killifish, fisherman's_lure"
$ws.Range("G14").Value = "raw_wound, adhesive_bandage, bandage, band_aid, oblique_bandage, flesh_wound, roller_bandage, laceration, capeline_bandage, elastoplast"
$ws.Range("G15").Value = "down, air_mass, high, solar_gravity, above, rise, lowering, rock_bottom, lower_berth, gravitation"
$ws.Range("G16").Value = "fast, quick, swift, speed, western_fence_lizard"
$ws.Range("G17").Value = "pitched_battle, battle, quarrel, fight, dispute, scuffle, melee, combat, controversy, fistfight"
$ws.Range("G18").Value = "bird, dickeybird, tailed_frog, uropygium, piciform_bird, protoavis, sinornis, milvus, tongueless_frog, caprimulgiform_bird"
$ws.Range("G19").Value = "thumb, fingertip, baseball_glove, thenar, nipa_palm, fishtail_palm, batting_glove, hand, lady_palm, royal_palm"
$ws.Range("G20").Value = "sleeping_pill, sleigh_bed, platform_bed, glutethimide, lake_bed, turnip_bed, sleeping, murphy_bed, ethchlorvynol, sleep"
$ws.Range("G21").Value = "talk, chitchat, chat, telephone_conversation, gossiping, table_talk, shmooze, shop_talk, scandalmonger, conversation"
$ws.Range("G22").Value = "fang, toxicognath, white_wolf, carnassial_tooth, timber_wolf, tooth, gingiva, malposed_tooth, red_wolf, dentition"
$ws.Range("G23").Value = "tobacco_juice, swamp, gleet, mud, mucus, wetland, bog, slick, spit, salt_marsh"
$ws.Range("G24").Value = "wordbook, pocket_dictionary, desk_dictionary, etymological_dictionary, bilingual_dictionary, thesaurus, learner's_dictionary, internal_rhyme, lexis, eye_rhyme"
$ws.Range("G25").Value = "wrong, error, mistake, erroneousness, unfairness, incorrectness, injustice"
$ws.Range("G26").Value = "tympanic_vein, sublingual_vein, stylomastoid_vein, subclavian_vein, brachial_vein, maxillary_vein, popliteal_vein, basal_vein, vertebral_vein, intercostal_vein"
$ws.Range("G27").Value = "crown_princess, princess_royal, empress, prince, emperor, second_empire, queen_consort, queen, prince_of_wales, grand_duchess"
$ws.Range("G28").Value = "music_stool, flat_bench, chaise_longue, settee, morris_chair, couch, banquette, footstool, recliner, campstool"
$ws.Range("G29").Value = "chemistry_lab, erlenmeyer_flask, vacuum_flask, bottle, hipflask, pill_bottle, specimen_bottle, flagon, drinking_vessel, phial"
$ws.Range("G31").Value = "grave, churchyard, diabetic_coma, burial_chamber, burial, coffin, semicoma, morgue, potter's_field, funeral"
$ws.Range("G32").Value = "reign_of_terror, fear, panic, horror, terrorization, intimidation, cold_feet, anxiety_attack"
$ws.Range("G33").Value = "finger, intercapitular_vein, oligodactyly, digital_arteries, big_toe, thumb, forepaw, fingertip, little_toe, hands"
$ws.Range("G34").Value = "lay_reader, side_chapel, lady_chapel, religious_doctrine, old_catholic_church, sacerdotalism, divine_messenger, old_catholic, catholic_church, church_service"
$ws.Range("G35").Value = "death's_head, human_head, torso, live_body, cranium, human_body, axial_skeleton, sinciput, musculoskeletal_system, coronal_suture"
$ws.Range("G36").Value = "violin, clarinet, wind_instrument, trombone, heckelphone, bowed_stringed_instrument, keyboard_instrument, musical_instrument, oboe, stradavarius"
$ws.Range("G37").Value = "paper_fastener, message_pad, pencil_sharpener, copyholder, staple_gun, pen, office_furniture, rolodex, writing_implement, notepad"
$ws.Range("G38").Value = "officer, sheriff, constable, police, police_sergeant, policeman, insignia, police_commissioner, military_officer, desk_sergeant"
$ws.Range("G39").Value = "delta_ray, free_electron, moment_of_inertia, momentum, angular_momentum"
$ws.Range("G40").Value = "perspiration, balanced_diet, carbohydrate_loading, allergy_diet, reducing_diet, stress, bland_diet, light_diet, dietary, overstrain"
$ws.Range("G41").Value = "homicide, contract_killing, aggravated_assault, crime, manslaughter, robbery, mugging, mariticide, killing, armed_robbery"
$ws.Range("G42").Value = "dentist's_drill, breast_drill, bore_bit, straight_flute, core_drill, power_drill, dig, drilling, posthole, burial_chamber"
$ws.Range("G43").Value = "great_care, scheduled_maintenance, tree_surgery, camera_care, kid_glove, tender_loving_care, personal_care, due_care, nurturance, carefulness"
$ws.Range("G44").Value = "white_wolf, moon, timber_wolf, red_wolf, cosmic_time"
$ws.Range("G45").Value = "chance, brass_ring, blossoming, blooming, flower, umbrellawort, occasion, divide, flourish, break"
$ws.Range("G46").Value = "cunning, common_dolphin, delphinus, river_dolphin"
$ws.Range("G47").Value = "clupeid_fish, pilchard, food_fish, sild, fish, scombroid, chum_salmon, black_duck, whitebait, saltwater_fish"
$ws.Range("G48").Value = "european_country, baltic_state, scandinavian_country, balkan_country, tartary, european, european_russia, stuffed_mushroom, europa, continent"
$ws.Range("G49").Value = "austenitic_steel, chisel_steel, medium_steel, brave, quenched_steel, fighter, mild_steel, combatant, tool_steel, crucible_steel"

# Update accuracy cell (H49) as text, not a parsed percentage number
$ws.Range("H49").Formula = "'39.6%"
$ws.Range("H49").Style = "Normal"

# Update selected cell to reflect the saved view
$ws.Range("G2").Select() | Out-Null

# Approximate the updated column widths for F and G
$ws.Columns("F:F").ColumnWidth = 10.83
$ws.Columns("G:G").ColumnWidth = 105.17
